$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 11
